# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column
# (value "stock" for the existing row) between the existing "total"
# and "date" columns. Inserting a whole column preserves the existing
# header/data cell styles (bold-bordered header style, plain data
# style) and shifts date / legislator_name / legislator_id one column
# to the right without disturbing their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H; this shifts the old H:J (date, legislator_name,
# legislator_id) to I:K while keeping each cell's original formatting.
$ws.Columns("H").Insert()

# Populate the newly inserted column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
